# Update Disease Ontology source_version to the August 2024 release
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "v2024-08-29"

# Move the active selection to E3 (as in the saved workbook)
$ws.Range("E3").Select()
